# "finish dev of gacha" — add casino-token-cost rows for normal/advanced
# gacha pulls to the intInit sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intInit")

# Values are entered in the same order the original author typed them
# (label, then description, for each row) so new shared strings land at
# the same indices as in the source edit.
$ws.Range("A13").Value = "casinoTokenNeededPerNormalGacha"
$ws.Range("C13").Value = "普通gacha一次需要多少赌币"
$ws.Range("C14").Value = "高级gacha一次需要多少赌币"
$ws.Range("A14").Value = "casinoTokenNeededPerAdvancedGacha"

$ws.Range("B13").Value = 100
$ws.Range("B14").Value = 1000

# Match the thin-border / centered formatting used by the rest of the table.
$ws.Range("A13:C14").Borders.LineStyle = 1
$ws.Range("A13:C14").HorizontalAlignment = -4108
$ws.Range("A13:C14").VerticalAlignment = -4108

$ws.Rows.Item(13).RowHeight = 20
$ws.Rows.Item(14).RowHeight = 20

$ws.Range("D13").Select()
